$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-27 08:18:43'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '67%'
$ws.Range('K2').Value = '0.3 MJ/m2'
$ws.Range('M2').Value = '4.0 °C 7:53 TU'
$ws.Range('O2').Value = '1.6 °C'
$ws.Range('E3').Value = '2026-02-27 08:18:46'
$ws.Range('I3').Value = '0.1 mm'
$ws.Range('K3').Value = '0.2 MJ/m2'
$ws.Range('M3').Value = '4.5 °C 7:59 TU'
$ws.Range('E4').Value = '2026-02-27 08:18:49'
$ws.Range('K4').Value = '0.1 MJ/m2'
$ws.Range('E5').Value = '2026-02-27 08:18:51'
$ws.Range('K5').Value = '0.1 MJ/m2'
$ws.Range('M5').Value = '4.7 °C 7:39 TU'
$ws.Range('O5').Value = '3.5 °C'
$ws.Range('E6').Value = '2026-02-27 08:18:54'
$ws.Range('K6').Value = '0.1 MJ/m2'
$ws.Range('M6').Value = '9.5 °C 7:53 TU'
$ws.Range('E7').Value = '2026-02-27 08:18:57'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '99%'
$ws.Range('K7').Value = '0.5 MJ/m2'
$ws.Range('E8').Value = '2026-02-27 08:18:59'
$ws.Range('K8').Value = '0.5 MJ/m2'
$ws.Range('O8').Value = '11.8 °C'
$ws.Range('E9').Value = '2026-02-27 08:19:02'
$ws.Range('M9').Value = '10.3 °C 7:59 TU'
$ws.Range('O9').Value = '8.8 °C'
$ws.Range('E10').Value = '2026-02-27 08:19:04'
$ws.Range('K10').Value = '0.2 MJ/m2'
$ws.Range('M10').Value = '9.8 °C 7:59 TU'
$ws.Range('O10').Value = '8.9 °C'
$ws.Range('E11').Value = '2026-02-27 08:19:07'
$ws.Range('E12').Value = '2026-02-27 08:19:09'
$ws.Range('M12').Value = '10.3 °C 7:58 TU'
$ws.Range('O12').Value = '8.7 °C'
$ws.Range('E13').Value = '2026-02-27 08:19:12'
$ws.Range('J13').Value = '1032.3 hPa'
$ws.Range('K13').Value = '0.5 MJ/m2'
$ws.Range('L13').Value = '6.5 km/h - 148º 7:43 TU'
$ws.Range('O13').Value = '-1.8 °C'
$ws.Range('E14').Value = '2026-02-27 08:19:15'
$ws.Range('E15').Value = '2026-02-27 08:19:18'
$ws.Range('M15').Value = '10.0 °C 7:55 TU'
$ws.Range('E16').Value = '2026-02-27 08:19:20'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '28%'
$ws.Range('E17').Value = '2026-02-27 08:19:23'
$ws.Range('K17').Value = '0.6 MJ/m2'
$ws.Range('E18').Value = '2026-02-27 08:19:26'
$ws.Range('K18').Value = '0.1 MJ/m2'
$ws.Range('E19').Value = '2026-02-27 08:19:28'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '81%'
$ws.Range('K19').Value = '0.2 MJ/m2'
$ws.Range('M19').Value = '9.6 °C 7:59 TU'
$ws.Range('E20').Value = '2026-02-27 08:19:31'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '55%'
$ws.Range('K20').Value = '0.5 MJ/m2'
$ws.Range('M20').Value = '3.9 °C 7:55 TU'
$ws.Range('O20').Value = '2.2 °C'
$ws.Range('E21').Value = '2026-02-27 08:19:33'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '80%'
$ws.Range('K21').Value = '0.4 MJ/m2'
$ws.Range('O21').Value = '3.5 °C'
$ws.Range('E22').Value = '2026-02-27 08:19:36'
$ws.Range('K22').Value = '0.4 MJ/m2'
$ws.Range('E23').Value = '2026-02-27 08:19:39'
$ws.Range('E24').Value = '2026-02-27 08:19:41'
$ws.Range('J24').Value = '1026.1 hPa'
$ws.Range('K24').Value = '0.4 MJ/m2'
$ws.Range('O24').Value = '3.9 °C'
$ws.Range('E25').Value = '2026-02-27 08:19:44'
$ws.Range('K25').Value = '0.5 MJ/m2'
$ws.Range('M25').Value = '8.0 °C 7:59 TU'
$ws.Range('O25').Value = '4.8 °C'
$ws.Range('E26').Value = '2026-02-27 08:19:47'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '41%'
$ws.Range('J26').Value = '1024.2 hPa'
$ws.Range('K26').Value = '0.4 MJ/m2'
$ws.Range('M26').Value = '12.7 °C 7:48 TU'
$ws.Range('O26').Value = '7.9 °C'
$ws.Range('E27').Value = '2026-02-27 08:19:50'
$ws.Range('K27').Value = '0.4 MJ/m2'
$ws.Range('M27').Value = '5.8 °C 7:54 TU'
$ws.Range('O27').Value = '4.3 °C'
$ws.Range('E28').Value = '2026-02-27 08:19:52'
$ws.Range('O28').Value = '5.2 °C'
$ws.Range('E29').Value = '2026-02-27 08:19:55'
$ws.Range('E30').Value = '2026-02-27 08:19:58'
$ws.Range('E31').Value = '2026-02-27 08:20:01'
$ws.Range('O31').Value = '9.3 °C'
$ws.Range('E32').Value = '2026-02-27 08:20:03'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '92%'
$ws.Range('K32').Value = '0.4 MJ/m2'
$ws.Range('M32').Value = '5.9 °C 7:48 TU'
$ws.Range('O32').Value = '0.8 °C'
$ws.Range('E33').Value = '2026-02-27 08:20:06'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '71%'
$ws.Range('J33').Value = '1029.5 hPa'
$ws.Range('K33').Value = '0.5 MJ/m2'
$ws.Range('M33').Value = '5.2 °C 7:46 TU'
$ws.Range('O33').Value = '2.1 °C'
$ws.Range('E34').Value = '2026-02-27 08:20:08'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '42%'
$ws.Range('K34').Value = '0.1 MJ/m2'
$ws.Range('L34').Value = '22.0 km/h - 26º 7:37 TU'
$ws.Range('O34').Value = '2.5 °C'
$ws.Range('E35').Value = '2026-02-27 08:20:11'
$ws.Range('H35').NumberFormat = '@'
$ws.Range('H35').Value = '40%'
$ws.Range('J35').Value = '1025.2 hPa'
$ws.Range('K35').Value = '0.2 MJ/m2'
$ws.Range('M35').Value = '12.1 °C 7:59 TU'
$ws.Range('O35').Value = '9.9 °C'
$ws.Range('E36').Value = '2026-02-27 08:20:14'
$ws.Range('K36').Value = '0.1 MJ/m2'
$ws.Range('M36').Value = '11.0 °C 7:38 TU'
$ws.Range('O36').Value = '9.6 °C'
$ws.Range('E37').Value = '2026-02-27 08:20:17'
$ws.Range('E38').Value = '2026-02-27 08:20:19'
$ws.Range('K38').Value = '0.1 MJ/m2'
$ws.Range('E39').Value = '2026-02-27 08:20:22'
$ws.Range('K39').Value = '0.6 MJ/m2'
$ws.Range('E40').Value = '2026-02-27 08:20:25'
$ws.Range('E41').Value = '2026-02-27 08:20:27'
$ws.Range('K41').Value = '0.2 MJ/m2'
$ws.Range('O41').Value = '8.3 °C'
$ws.Range('E42').Value = '2026-02-27 08:20:30'
$ws.Range('M42').Value = '10.7 °C 7:55 TU'
$ws.Range('O42').Value = '8.7 °C'
$ws.Range('E43').Value = '2026-02-27 08:20:33'
$ws.Range('H43').NumberFormat = '@'
$ws.Range('H43').Value = '99%'
$ws.Range('K43').Value = '0.4 MJ/m2'
$ws.Range('M43').Value = '6.0 °C 7:52 TU'
$ws.Range('O43').Value = '3.6 °C'
$ws.Range('E44').Value = '2026-02-27 08:20:36'
$ws.Range('K44').Value = '0.1 MJ/m2'
$ws.Range('E45').Value = '2026-02-27 08:20:38'
$ws.Range('E46').Value = '2026-02-27 08:20:41'
$ws.Range('K46').Value = '0.1 MJ/m2'
